$wb = $excel.ActiveWorkbook

$wsCred = $wb.Worksheets.Item("Credentials")
$wsCommon = $wb.Worksheets.Item("Common")
$wsUtil = $wb.Worksheets.Item("BEN_FTR_001_Utilities")

# ---------------------------------------------------------------
# "Common" sheet: insert a new Helix_Test_Case_Number column (D),
# shifting the old D:G columns to E:H, and turn the existing BEN id
# column (B) into the Helix-style shared-string labels.
# ---------------------------------------------------------------
$wsCommon.Columns.Item(4).Insert()
$wsCommon.Columns.Item(4).ColumnWidth = 23.25

$wsCommon.Range("D1").Value = "Helix_Test_Case_Number"
$wsCommon.Range("D2").Value = 22746
$wsCommon.Range("D3").Value = 1221
$wsCommon.Range("D4").Value = 4321

# D3 has no pre-existing neighbour cell to inherit a style from (row 3
# has no C3), so pull the numeric-cell style from C2 explicitly.
$wsCommon.Range("C2").Copy()
$wsCommon.Range("D3").PasteSpecial(-4122)

$wsCommon.Range("B2").Value = "BEN_22746"
$wsCommon.Range("B3").Value = "BEN_1221"
$wsCommon.Range("B4").Value = "BEN_4321"

# ---------------------------------------------------------------
# "BEN_FTR_001_Utilities" sheet: same new column + BEN-id relabelling.
# ---------------------------------------------------------------
$wsUtil.Columns.Item(4).Insert()
$wsUtil.Columns.Item(4).ColumnWidth = 38.42

$wsUtil.Range("D1").Value = "Helix_Test_Case_Number"
$wsUtil.Range("D2").Value = 22746
$wsUtil.Range("D3").Value = 1221
$wsUtil.Range("D4").Value = 4321

$wsUtil.Range("B2").Value = "BEN_22746"
$wsUtil.Range("B3").Value = "BEN_1221"
$wsUtil.Range("B4").Value = "BEN_4321"

# ---------------------------------------------------------------
# Selection/view bookkeeping. Selecting a range on a sheet activates
# it, so touch the non-active sheets first and finish on the sheet
# that should stay tabSelected ("Credentials").
# ---------------------------------------------------------------
$wsUtil.Range("D2:D4").Select()
$wsCommon.Range("B2").Select()
$wsCred.Range("C1:C1048576").Select()
